$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin'
$ws.Range("D2").Value = '91.040.26'

# Row 3: 'Ethereum'
$ws.Range("D3").Value = '3.188.54'
$ws.Range("E3").Value = '  -2.86%  '

# Row 4: 'TetherUSD'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.17%  '

# Row 5: 'Solana'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.13%  '

# Row 6: 'BNB'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '617.75'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.51%  '

# Row 7: 'Dogecoin'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.390'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.27%  '

# Row 8: 'XRP'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.690'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.38%  '

# Row 9: 'USDC'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.998'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.02%  '

# Row 10: 'LidoStakedEther'
$ws.Range("D10").Value = '3.180.98'
$ws.Range("E10").Value = '  -2.96%  '

# Row 11: 'Cardano'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.579'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.21%  '

# Row 12: 'TRON'
$ws.Range("E12").Value = '  -5.01%  '

# Row 13: 'ShibaInu'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000253'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.70%  '

# Row 14: 'WrappedBTC'
$ws.Range("D14").Value = '90.506.80'
$ws.Range("E14").Value = '  +1.48%  '

# Row 15: 'WrappedliquidstakedEther2.0'
$ws.Range("D15").Value = '3.761.08'
$ws.Range("E15").Value = '  -3.13%  '

# Row 16: 'Avalanche'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '32.89'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.74%  '

# Row 17: 'Toncoin'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.24'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.59%  '

# Row 18: 'WrappedEther'
$ws.Range("D18").Value = '3.161.59'
$ws.Range("E18").Value = '  -4.12%  '

# Row 19: 'SuiNetwork'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.26'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.42%  '

# Row 20: 'Chainlink'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.50'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.23%  '

# Row 21: 'BitcoinCash'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '441.17'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.90%  '

# Row 22: 'PEPE'
$ws.Range("E22").Value = '  +38.93%  '

# Row 23: 'Uniswap'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.59'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.50%  '

# Row 24: 'Polkadot'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.11'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.47%  '

# Row 25: 'NEARProtocol'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.15'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.19%  '

# Row 26: 'Aptos'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.77'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.53%  '

# Row 27: 'WrappedeETH'
$ws.Range("D27").Value = '3.347.73'
$ws.Range("E27").Value = '  -3.37%  '

# Row 28: 'Litecoin'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '74.78'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.81%  '

# Row 29: 'Dai'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.06%  '

# Row 30: 'Cronos'
$ws.Range("E30").Value = '  -6.04%  '

# Row 31: 'Binance-PegBSC-USD'
$ws.Range("E31").Value = '  +0.03%  '

# Row 32: 'dogwifhat'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.17'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +35.24%  '

# Row 33: 'InternetComputer(DFINITY)'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.46'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.51%  '

# Row 34: 'Bittensor'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '533.49'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.10%  '

# Row 35: 'RenderToken'
$ws.Range("E35").Value = '  -2.20%  '

# Row 36: 'PancakeSwap'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.87'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.83%  '

# Row 37: 'Fetch.AI'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.25'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -8.76%  '

# Row 38: 'WhiteBITCoin' -> 'EthereumClassic'
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '21.96'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.32%  '

# Row 39: 'EthereumClassic' -> 'WhiteBITCoin'
$ws.Range("B39").Value = 'WhiteBITCoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.32'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.24%  '

# Row 40: 'FirstDigitalUSD' -> 'Kaspa'
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.128'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -8.53%  '

# Row 41: 'Kaspa' -> 'FirstDigitalUSD'
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.997'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.20%  '

# Row 42: 'USDe'
$ws.Range("E42").Value = '  -0.10%  '

# Row 43: 'PolygonEcosystemToken'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.376'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.29%  '

# Row 44: 'Stacks'
$ws.Range("E44").Value = '  -5.74%  '

# Row 45: 'Monero'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '146.77'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.38%  '

# Row 46: 'OKB'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '44.38'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.37%  '

# Row 47: 'Aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '172.69'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.36%  '

# Row 48: 'Stellar'
$ws.Range("E48").Value = '  -9.17%  '

# Row 49: 'ImmutableX'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.24'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.35%  '

# Row 50: 'ARBITRUM'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.613'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.03%  '

# Row 51: 'Filecoin'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.08'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.86%  '
